# Update the ms_refs table:
#  - add "US Core Device Profile" row (profile -> US Core Patient Profile)
#  - add "US Core FamilyMemberHistory Profile" row (profile -> US Core Patient Profile)
#  - remove the "US Core Implantable Device Profile" row (duplicate/obsolete)
#  - add "US Core PMO ServiceRequest Profile" row with its target references
#  - fix "US Core Provenance Profile" targets to include the "Resource" reference
#  - renumber the leading index column to stay sequential

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert "US Core Device Profile" before row 19 (DiagnosticReport ... Laboratory Results Reporting)
$ws.Cells.Item(19, 1).EntireRow.Insert()
$ws.Cells.Item(19, 2).Value = "US Core Device Profile"
$ws.Cells.Item(19, 3).Value = "US Core Patient Profile"

# 2) Insert "US Core FamilyMemberHistory Profile" before row 24 (US Core Goal Profile)
$ws.Cells.Item(24, 1).EntireRow.Insert()
$ws.Cells.Item(24, 2).Value = "US Core FamilyMemberHistory Profile"
$ws.Cells.Item(24, 3).Value = "US Core Patient Profile"

# 3) Delete the obsolete "US Core Implantable Device Profile" row (now row 29)
$ws.Cells.Item(29, 1).EntireRow.Delete()

# 4) Insert "US Core PMO ServiceRequest Profile" before row 43 (US Core Practitioner Profile)
$ws.Cells.Item(43, 1).EntireRow.Insert()
$ws.Cells.Item(43, 2).Value = "US Core PMO ServiceRequest Profile"
$ws.Cells.Item(43, 3).Value = "US Core ADI DocumentReference Profile"
$ws.Cells.Item(43, 4).Value = "US Core Encounter Profile"
$ws.Cells.Item(43, 5).Value = "US Core Patient Profile"
$ws.Cells.Item(43, 6).Value = "US Core Practitioner Profile"

# 5) Fix "US Core Provenance Profile" row (row 47): insert "Resource" at target1,
#    shift the previous "US Core Organization Profile" value to target2
$ws.Cells.Item(47, 3).Value = "Resource"
$ws.Cells.Item(47, 4).Value = "US Core Organization Profile"

# 6) Renumber the leading index column (A) sequentially starting at 0 for every data row
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
